# Project Log update: subtitle-sync fix follow-up + two new log rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: Matthew's hours corrected 4 -> 6 -----------------------------
$ws.Range("D31").Value = 6

# --- Row 34: new "Development" entry (26/6/2019) --------------------------
$ws.Range("B34").Value = "Development"

$ws.Range("C30").Copy()
$ws.Range("C34").PasteSpecial(-4122)   # xlPasteFormats - reuse the date style
$ws.Range("C34").Value = 43642         # 26/6/2019

$ws.Range("D34").Value = 6

$ws.Range("F30").Copy()
$ws.Range("F34").PasteSpecial(-4122)   # xlPasteFormats - reuse the wrap-text style
$ws.Range("F34").Value = "1) Bug fixes across the application with regards to delay-disabler.js and how the application detects if a webcast is single or double video stream etc`n2) In-depth testing of silence removal detects memory leak`n3) Further testing and research reveals that the leak is caused by OfflineAudioContext not being garbage collected"

$ws.Rows.Item(34).RowHeight = 72.5

# --- Row 35: new "Development" entry (27/6/2019 - 28/6/2019) --------------
$ws.Range("B35").Value = "Development"

$ws.Range("C35").HorizontalAlignment = -4108   # xlCenter
$ws.Range("C35").VerticalAlignment = -4108     # xlCenter
$ws.Range("C35").WrapText = $true
$ws.Range("C35").NumberFormat = "mm-dd-yy"     # canonical code for built-in numFmtId 14
$ws.Range("C35").Value = "27/6/2019 - `n28/6/2019"

$ws.Range("D35").Value = 14

$ws.Range("F30").Copy()
$ws.Range("F35").PasteSpecial(-4122)   # xlPasteFormats - reuse the wrap-text style
$ws.Range("F35").Value = "1) Further developmental work to remove memory leak; re-implemented using a buffering system to reduce amount of threads spawned, but memory leak persisted.`n2) Re-implemented using a web worker instead and analysing the PCM data directly, bypassing the OfflineAudioContext, which fixed memory leak. This had the side effect of making the code more readable."

$ws.Rows.Item(35).RowHeight = 87

# --- Totals (D43/C43 formulas recalc automatically, but nudge the view) ---
$ws.Range("D43").Formula = "=SUM(D3:D42)"
$ws.Range("C43").Formula = "=SUM(D43:E43)"

# --- Sheet view: scrolled one column right, selection moved to F36 --------
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("F36").Select()
